$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B39 to be a true numeric value (3) instead of text "3"
$ws.Range("B39").Value = 3

# Insert new row 40 with the data that was previously in row 39 for column B (text "3"),
# plus the new annotation data.
$ws.Range("A40").Value = "Ruilin"
# Force B40 to stay text "3" (like the legacy B39 cell) instead of being
# auto-coerced to a number, then strip the temporary text number-format
# back off so no extra style is left behind on the cell.
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "3"
$ws.Range("B40").ClearFormats()
$ws.Range("C40").Value = "无"
$ws.Range("D40").Value = "DIS"
$ws.Range("E40").Value = "MET"
$ws.Range("F40").Value = "42b1e2ab-785d-481e-b197-1cf6913a8b3e"
$ws.Range("G40").Value = "SJQO7UJCW_annotated.xlsx"
$ws.Range("H40").Value = "However, our main point of the paper is to demonstrate the effectiveness of proposed method against our baseline model shown in Table 1 and 2."
